$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Portronics - Wireless Bluetooth Headset"
$ws.Range("B1").Value = "Rs. 1,149"

$ws.Range("A2").Value = "boAt Airdopes 131/138 Twin Wireless Earbuds with IWP Technology, Bluetooth V5.0, Immersive Audio, Up to 15H Total Playback, Instant Voice Assistant and Type-C Charging,Bluetooth Earphone (Active Black)"
$ws.Range("B2").Value = "Rs. 1,499"

$ws.Range("A3").Value = "Tecsox PowerHouse Earbud In Ear Bluetooth Earphone 45 Hours Playback Bluetooth IPX5(Splash Proof) Powerfull Bass -Bluetooth V 5.1 Black"
$ws.Range("B3").Value = "Rs. 725"

$ws.Range("A4").Value = "hitage NBT-6586+ Neckband In Ear Bluetooth Neckband 22 Hours Playback IPX6(Water Resistant) Fast charging -Bluetooth V 5.0 Red"
$ws.Range("B4").Value = "Rs. 789"

$ws.Range("A5").Value = "pTron Bassbuds Duo Ear Buds Wireless With Mic Headphones/Earphones White"
$ws.Range("B5").Value = "Rs. 799"
